# Commit: "Include abby error rate in resutls"
#
# Adds a new "abby_data" worksheet at the end of the workbook (after
# "error" and "extra_work"), containing the abby error-rate numbers laid
# out the same way as the other two result sheets: a bold/boxed header
# row of "<n>lpp" labels in B1:F1, a boxed "0" in A2, and the rate values
# in B2:F2.

$wb = $excel.ActiveWorkbook

$errorSheet = $wb.Worksheets.Item(1)

# Add the new sheet after the last existing sheet ("extra_work") so it
# lands at the end of the tab strip, matching the diff ordering.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "abby_data"

# Reuse the existing bold / bordered / centered look from the "error"
# sheet's header row and leading A-column cell instead of re-declaring
# the formatting from scratch (keeps the same style entry, no new styles
# added to the workbook).
$errorSheet.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$errorSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)     # xlPasteFormats

# Header row
$ws.Range("B1").Value = "14lpp"
$ws.Range("C1").Value = "151lpp"
$ws.Range("D1").Value = "397lpp"
$ws.Range("E1").Value = "665lpp"
$ws.Range("F1").Value = "996lpp"

# Data row
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.256
$ws.Range("C2").Value = 0.1
$ws.Range("D2").Value = 0.1230769230769231
$ws.Range("E2").Value = 0.1359447004608295
$ws.Range("F2").Value = 0.08415841584158416
